$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need NumberFormat forced to
# Text ("@") first, otherwise Excel COM auto-converts the assigned string into a
# numeric value instead of keeping it as text (matching the source data which
# stores these as strings).
$ws.Range("D2").Value = '42.826.01'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '2.542.03'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.92'
$ws.Range("E5").Value = '  +1.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.74'
$ws.Range("E6").Value = '  +6.44%  '
$ws.Range("E7").Value = '  +0.57%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.545'
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.79'
$ws.Range("E10").Value = '  +1.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0828'
$ws.Range("E11").Value = '  +3.54%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.115'
$ws.Range("E12").Value = '  +2.14%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.62'
$ws.Range("E13").Value = '  -0.92%  '
$ws.Range("D14").Value = '2.933.87'
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").Value = '2.553.68'
$ws.Range("E15").Value = '  +1.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.06'
$ws.Range("E16").Value = '  +6.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.866'
$ws.Range("E17").Value = '  -0.68%  '
$ws.Range("D18").Value = '42.867.15'
$ws.Range("E18").Value = '  +0.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.58'
$ws.Range("E19").Value = '  +4.49%  '
$ws.Range("D20").Value = '0.0₃0992'
$ws.Range("E20").Value = '  +1.19%  '
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.00'
$ws.Range("E22").Value = '  +0.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.91'
$ws.Range("E23").Value = '  -0.32%  '
$ws.Range("E24").Value = '  +1.68%  '
$ws.Range("E25").Value = '  -1.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.12'
$ws.Range("E26").Value = '  -3.69%  '
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.19'
$ws.Range("E28").Value = '  +1.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.86'
$ws.Range("E29").Value = '  +2.76%  '
$ws.Range("E30").Value = '  -5.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.14'
$ws.Range("E31").Value = '  +3.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '157.92'
$ws.Range("E32").Value = '  +3.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.47'
$ws.Range("E33").Value = '  +14.80%  '
$ws.Range("E34").Value = '  -2.36%  '
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.30'
$ws.Range("E36").Value = '  -2.58%  '
$ws.Range("E37").Value = '  -4.48%  '
$ws.Range("E38").Value = '  +1.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '25.53'
$ws.Range("E39").Value = '  +6.97%  '
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("E42").Value = '  +29.53%  '
$ws.Range("E43").Value = '  +0.72%  '
$ws.Range("D44").Value = '2.095.88'
$ws.Range("E44").Value = '  +0.77%  '
$ws.Range("E45").Value = '  -1.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.28'
$ws.Range("E47").Value = '  +4.68%  '
$ws.Range("E48").Value = '  -0.81%  '
$ws.Range("D49").Value = '2.792.25'
$ws.Range("E49").Value = '  +0.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.65'
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.68'
$ws.Range("E51").Value = '  +1.36%  '
